$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 updates
$ws.Range("F3").Value = 1.73
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 3.35
$ws.Range("K3").Value = 4.1
$ws.Range("P3").Value = 1.72
$ws.Range("Q3").Value = 2.12

# Row 4 updates
$ws.Range("F4").Value = 1.76
$ws.Range("G4").Value = 1.96
$ws.Range("H4").Value = 4.6
$ws.Range("I4").Value = 7
$ws.Range("K4").Value = 4.3
